# Updated crypto price/volume figures per latest scrape run.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Some 'Price' cells hold digit-grouped text (e.g. "1.003") that Excel's
# automatic type inference would otherwise coerce into a number when
# assigned through .Value. Temporarily force those specific cells to the
# Text number format so the string is preserved verbatim, then restore
# the cell style to Normal once the value is in place.
$textPriceCells = @("D5", "D6", "D7", "D8", "D9", "D10", "D13", "D14", "D15", "D18", "D19", "D21", "D22", "D23", "D26", "D27", "D28", "D29", "D30", "D32", "D33", "D34", "D35", "D36", "D37", "D38", "D39", "D40", "D41", "D42", "D43", "D44", "D46", "D47", "D49", "D50", "D51")
foreach ($addr in $textPriceCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = '22.411.93'
$ws.Range("E2").Value = '  +0.24%  '
$ws.Range("D3").Value = '1.570.95'
$ws.Range("E4").Value = '  +0.20%  '
$ws.Range("D5").Value = '1.003'
$ws.Range("E5").Value = '  +0.22%  '
$ws.Range("D6").Value = '290.21'
$ws.Range("E6").Value = '  -0.01%  '
$ws.Range("D7").Value = '0.3746'
$ws.Range("E7").Value = '  -0.76%  '
$ws.Range("D8").Value = '49.62'
$ws.Range("E8").Value = '  +0.23%  '
$ws.Range("D9").Value = '0.3372'
$ws.Range("E9").Value = '  -1.65%  '
$ws.Range("D10").Value = '0.07482'
$ws.Range("E10").Value = '  -2.68%  '
$ws.Range("E11").Value = '  -3.34%  '
$ws.Range("E12").Value = '  +0.18%  '
$ws.Range("D13").Value = '21.03'
$ws.Range("E13").Value = '  -1.80%  '
$ws.Range("D14").Value = '5.928'
$ws.Range("E14").Value = '  -1.80%  '
$ws.Range("D15").Value = '6.874'
$ws.Range("E15").Value = '  -1.10%  '
$ws.Range("D16").Value = '1.569.58'
$ws.Range("E16").Value = '  -0.29%  '
$ws.Range("E17").Value = '  -1.98%  '
$ws.Range("D18").Value = '89.44'
$ws.Range("E18").Value = '  -1.30%  '
$ws.Range("D19").Value = '0.06690'
$ws.Range("E19").Value = '  -0.45%  '
$ws.Range("D21").Value = '6.172'
$ws.Range("E21").Value = '  -1.60%  '
$ws.Range("D22").Value = '16.18'
$ws.Range("E22").Value = '  -3.02%  '
$ws.Range("D23").Value = '11.86'
$ws.Range("E23").Value = '  -1.16%  '
$ws.Range("D24").Value = '22.415.27'
$ws.Range("E24").Value = '  +0.20%  '
$ws.Range("E25").Value = '  -1.02%  '
$ws.Range("D26").Value = '2.553'
$ws.Range("E26").Value = '  -9.26%  '
$ws.Range("D27").Value = '20.04'
$ws.Range("E27").Value = '  -1.62%  '
$ws.Range("D28").Value = '147.04'
$ws.Range("E28").Value = '  +1.49%  '
$ws.Range("D29").Value = '5.000'
$ws.Range("E29").Value = '  -0.65%  '
$ws.Range("D30").Value = '124.71'
$ws.Range("E30").Value = '  -1.23%  '
$ws.Range("D31").Value = '1.745.70'
$ws.Range("E31").Value = '  -0.15%  '
$ws.Range("D32").Value = '0.9923'
$ws.Range("E32").Value = '  -3.04%  '
$ws.Range("D33").Value = '1.956'
$ws.Range("E33").Value = '  -3.19%  '
$ws.Range("D34").Value = '5.931'
$ws.Range("E34").Value = '  -5.24%  '
$ws.Range("D35").Value = '9.748'
$ws.Range("E35").Value = '  -3.86%  '
$ws.Range("D36").Value = '0.08428'
$ws.Range("E36").Value = '  -2.01%  '
$ws.Range("D37").Value = '1.381'
$ws.Range("E37").Value = '  +4.54%  '
$ws.Range("D38").Value = '0.02451'
$ws.Range("E38").Value = '  -3.62%  '
$ws.Range("D39").Value = '0.06481'
$ws.Range("E39").Value = '  +1.36%  '
$ws.Range("D40").Value = '0.2253'
$ws.Range("E40").Value = '  -3.33%  '
$ws.Range("D41").Value = '5.395'
$ws.Range("E41").Value = '  -3.44%  '
$ws.Range("D42").Value = '11.33'
$ws.Range("E42").Value = '  -3.42%  '
$ws.Range("D43").Value = '0.6227'
$ws.Range("E43").Value = '  -2.86%  '
$ws.Range("D44").Value = '13.97'
$ws.Range("E44").Value = '  -2.06%  '
$ws.Range("E45").Value = '  +0.20%  '
$ws.Range("D46").Value = '3.813'
$ws.Range("E46").Value = '  +1.23%  '
$ws.Range("D47").Value = '0.5792'
$ws.Range("E47").Value = '  -3.69%  '
$ws.Range("D49").Value = '125.52'
$ws.Range("E49").Value = '  +0.88%  '
$ws.Range("D50").Value = '1.227'
$ws.Range("E50").Value = '  -7.43%  '
$ws.Range("D51").Value = '0.07303'
$ws.Range("E51").Value = '  +0.15%  '

foreach ($addr in $textPriceCells) {
    $ws.Range($addr).Style = "Normal"
}
